$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21529.55234596766
$ws.Range("B3").Value = 12119.65720836722
$ws.Range("B4").Value = 23560.23798410574
$ws.Range("B5").Value = 26269.89319920146
$ws.Range("B6").Value = 13825.05799014997
$ws.Range("B7").Value = 16003.72212012359
$ws.Range("B8").Value = 31591.08841752559
$ws.Range("B9").Value = 19312.6427917023
$ws.Range("B10").Value = 17417.28114820644
$ws.Range("B11").Value = 21565.59850193837
$ws.Range("B12").Value = 33159.04357116088
$ws.Range("B13").Value = 18235.27350696297
$ws.Range("B14").Value = 25812.53772514699
$ws.Range("B15").Value = 29309.17501649283
$ws.Range("B16").Value = 19619.86323274543
$ws.Range("B17").Value = 29537.67432737427
$ws.Range("B18").Value = 23059.07907315761
$ws.Range("B19").Value = 25586.8417362342
$ws.Range("B20").Value = 23278.0677016509
$ws.Range("B21").Value = 21995.90426069828
$ws.Range("B22").Value = 11662.51018416592
$ws.Range("B23").Value = 28153.61859544036
$ws.Range("B24").Value = 23475.9361988595
$ws.Range("B25").Value = 15973.33497959356
$ws.Range("B26").Value = 19355.75462398761
$ws.Range("B27").Value = 23969.13432540124
$ws.Range("B28").Value = 28011.25851730144
$ws.Range("B29").Value = 14127.31614997417
$ws.Range("B30").Value = 25211.77921435846
$ws.Range("B31").Value = 19501.35876990029
$ws.Range("B32").Value = 23939.44742569692
$ws.Range("B33").Value = 18850.66251960217
$ws.Range("B34").Value = 24233.21425236989
$ws.Range("B35").Value = 16126.70353177607
$ws.Range("B36").Value = 31386.58785067385
$ws.Range("B37").Value = 14693.84896910398
$ws.Range("B38").Value = 24012.79068072772
$ws.Range("B39").Value = 19968.87984878441
$ws.Range("B40").Value = 32745.67318259205
$ws.Range("B41").Value = 28710.69724942467
$ws.Range("B42").Value = 14686.44703493291
$ws.Range("B43").Value = 27029.98206282103
$ws.Range("B44").Value = 26573.22980937499
$ws.Range("B45").Value = 28743.77539523179
$ws.Range("B46").Value = 25805.84781824834
$ws.Range("B47").Value = 20032.92569838882
$ws.Range("B48").Value = 13189.20463995007
$ws.Range("B49").Value = 30238.49035583191
$ws.Range("B50").Value = 30700.60215244205
$ws.Range("B51").Value = 24084.25971978603
$ws.Range("B52").Value = 17223.53343880632
$ws.Range("B53").Value = 21014.40249785926
$ws.Range("B54").Value = 29172.69486243113
$ws.Range("B55").Value = 30142.86781587759
$ws.Range("B56").Value = 30359.81777166135
$ws.Range("B57").Value = 26049.83187053042
$ws.Range("B58").Value = 26312.51722238317
$ws.Range("B59").Value = 18757.97189773978
$ws.Range("B60").Value = 14116.17790631904
$ws.Range("B61").Value = 29921.75883401286
